$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto symbol table (price, 1h volume %, and hour columns) with the
# latest scrape, as produced by the scheduled GitHub Actions job. The sheet stores
# these figures as text (not numbers) so values like "--" coexist with formatted
# numbers/percentages in the same columns. Setting NumberFormat to "@" (Text) on
# each target cell before assigning its value keeps Excel from auto-converting
# numeric-looking strings into real numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.88%"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G2").Value = "4"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.77"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "0.64%"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G3").Value = "4"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.179"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.83%"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("G4").Value = "4"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07480"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.54%"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "4"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.421"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "40.88%"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("G6").Value = "4"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.019"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.05%"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("G7").Value = "4"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9141"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.46%"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("G8").Value = "4"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1734"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.66%"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "4"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07665"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.44%"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("G10").Value = "4"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08179"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "3.21%"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "4"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03038"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.11%"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "4"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09945"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.50%"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("G13").Value = "4"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001526"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2.20%"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("G14").Value = "4"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006072"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-2.80%"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("G15").Value = "4"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.503"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.34%"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "4"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.867"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.93%"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("G17").Value = "4"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.237"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.33%"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "4"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3262"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.90%"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("G19").Value = "4"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.31%"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("G20").Value = "4"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.656"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.76%"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("G21").Value = "4"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04606"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.86%"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("G22").Value = "4"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1566"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.10%"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("G23").Value = "4"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001261"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "3.86%"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("G24").Value = "4"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004536"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.59%"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("G25").Value = "4"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.15%"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("G26").Value = "4"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0002740"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "51.60%"
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "4"
$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "4"
$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "4"
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "4"
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "4"
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "4"
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "4"
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "4"
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "4"
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "4"
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "4"
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "4"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01762"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "5.47%"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "4"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04546"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-0.11%"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("G40").Value = "4"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007449"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "7.54%"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("G41").Value = "4"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1365"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.73%"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("G42").Value = "4"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002119"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.91%"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("G43").Value = "4"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01090"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-15.45%"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("G44").Value = "4"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006480"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "6.75%"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("G45").Value = "4"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-57.22%"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("G46").Value = "4"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("G47").Value = "4"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "4"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "4"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "4"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "4"
